# Update the "想去人数" (want-to-go count) figures that changed between
# the previous and newly generated scrape of the 南宁-漫展信息 data.
#
# Sheet "展览"   (sheet1): rows 4,6,7,10,11,12 column F
# Sheet "演出"   (sheet2): row 2 column F
# Sheet "全部类型" (sheet4): rows 4,6,7,9,11,12,13 column F
# Sheet "本地生活" (sheet3): unchanged

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 268
$wsExhibition.Range("F6").Value = 3179
$wsExhibition.Range("F7").Value = 2087
$wsExhibition.Range("F10").Value = 1182
$wsExhibition.Range("F11").Value = 214
$wsExhibition.Range("F12").Value = 1066

$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 32

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 268
$wsAll.Range("F6").Value = 3179
$wsAll.Range("F7").Value = 2087
$wsAll.Range("F9").Value = 32
$wsAll.Range("F11").Value = 1182
$wsAll.Range("F12").Value = 214
$wsAll.Range("F13").Value = 1066
